# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E) and "Correspond Handback DateTime" (H)
# timestamps for the c24ac456... row on the "zh-cn" sheet, and for the
# c24ac456... row on the "de-de" sheet, reflecting the newly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-20 17:07:50"
$wsZh.Range("H3").Value = "2016-03-20 17:08:31"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-20 17:07:58"
$wsDe.Range("H3").Value = "2016-03-20 17:08:44"
